$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Obras en general")

# --- Row 9: clear the "answered nothing" placeholder cells ---
# (these were empty inline-string cells in the source; the form resubmission
# that produced row 10 also cleaned these up on row 9)
$ws.Cells.Item(9, 16).ClearContents()   # P9
$ws.Cells.Item(9, 17).ClearContents()   # Q9
$ws.Cells.Item(9, 20).ClearContents()   # T9
$ws.Cells.Item(9, 21).ClearContents()   # U9
$ws.Cells.Item(9, 22).ClearContents()   # V9
$ws.Cells.Item(9, 23).ClearContents()   # W9

# --- Row 10: new record (WhatsApp del profesional saved) ---
function Set-TextCell($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextCell 10 1 "05/06/2025"                 # A10 Fecha
$ws.Cells.Item(10, 2).Value = "Ingeniero"       # B10 Profesión
$ws.Cells.Item(10, 3).Value = "Físico"          # C10 Formato
Set-TextCell 10 4 "5"                           # D10 Nro de Copias
$ws.Cells.Item(10, 5).Value = "Obra nueva"      # E10 Tipo de trabajo
$ws.Cells.Item(10, 6).Value = "VITALE JUAN ANTONIO"   # F10 Nombre del Profesional
$ws.Cells.Item(10, 7).Value = "NILOS ROBERTO"         # G10 Nombre del Comitente
$ws.Cells.Item(10, 8).Value = "FRANCISCO DE HARO 2745" # H10 Ubicación
$ws.Cells.Item(10, 9).Value = "1515/J/25"       # I10 Nro de expte municipal
Set-TextCell 10 11 "61518"                      # K10 Nro de partida inmobiliaria
Set-TextCell 10 12 "15000"                      # L10 Tasa de sellado
Set-TextCell 10 14 "15000"                      # N10 Visado de instalacion de Gas
Set-TextCell 10 15 "2000"                       # O10 Visado de instalacion de Salubridad
$ws.Cells.Item(10, 18).Value = "No pagado"      # R10 Estado pago sellado
$ws.Cells.Item(10, 19).Value = "No pagado"      # S10 Estado pago visado
Set-TextCell 10 25 "3764251817"                 # Y10 WhatsApp Profesional
